# "signed off time sheets"
#
# The Supervisor ("Prakruti Sinha") has now signed off on Roger Tan's
# week-4 timesheet:
#   - the Supervisor Name field (G6) is filled in
#   - the Supervisor Signature line (A27) is initialled "P.S"
#   - the Supervisor Signature date (D27) is filled in (12-Feb-2014)
#
# The new cells are formatted to match the existing, already-filled-in
# Employee signature block (A25 = name/signature style, D25 = date style),
# so copy those formats across before setting the values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly")

# --- Supervisor Name (G6), formatted like the Employee Name value (G4/A25) ---
$ws.Range("A25").Copy()
$ws.Range("G6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G6").Value2 = "Prakruti Sinha"

# --- Supervisor Signature (A27), formatted like the Employee Signature (A25) ---
$ws.Range("A25").Copy()
$ws.Range("A27").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A27").Value2 = "P.S"

# --- Supervisor Signature date (D27), formatted like the Employee date (D25) ---
$ws.Range("D25").Copy()
$ws.Range("D27").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D27").Value2 = 41682

# Clear clipboard/marching-ants state left behind by Copy()
$excel.CutCopyMode = 0

# Restore the view: scroll back to the top and leave the selection on M12,
# matching where the signer's cursor ended up after signing off.
$ws.Range("A1").Select()
$ws.Range("M12").Select()
